# "Agrego aproximacion a calendario academico"
#
# - Marks every currently-visible data row (condicion = "carrera", i.e. the
#   rows not hidden by the existing autofilter) with a 1 in column J.
# - Clears the autofilter criteria (keeps the filter range/dropdowns) and
#   unhides every row that the filter had hidden, except the trailing blank
#   row 162.
# - Moves the active selection to J159.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that were left visible by the old "condicion = carrera" filter -
# stamp a 1 into column J for each of them.
$flaggedRows = @(2,3,4,6,7,8,9,10,11,12,13,15,20,21,22,23,24,25,26,121,122,124,142,146,148,149,150,151,152,153,154,155,156,157,159)

foreach ($r in $flaggedRows) {
    $ws.Cells.Item($r, 10).Value = 1
}

# Drop the filter criteria (colId 4 -> "carrera") but leave the autofilter
# dropdown range A1:O161 in place.
$ws.ShowAllData()

# Unhide every data row hidden by the old filter, row 162 (the trailing
# blank row) stays hidden.
for ($r = 1; $r -le 161; $r++) {
    $ws.Rows.Item($r).Hidden = $false
}
$ws.Rows.Item(162).Hidden = $true

# Move the selection like the author left it.
$ws.Range("J159").Select()
